# Activity Log updates: arithmetic-unit work entries (29/03 - 31/03/2020)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# --- Row 6 (29/03/2020): "9:20pm"/"9:50pm" text times become real time-of-day
#     serial values (same displayed time, now numeric so it behaves as a time).
$ws.Range("D6").Value = 0.88888888888888884   # 9:20 PM
$ws.Range("E6").Value = 0.90972222222222221   # 9:50 PM

# --- Row 7 (29/03/2020): "9:55pm"/"10:45pm" text times -> numeric time values.
$ws.Range("D7").Value = 0.91319444444444453   # 9:55 PM
$ws.Range("E7").Value = 0.94791666666666663   # 10:45 PM

# --- Row 8 (new entry, 30/03/2020): full adder / ripple adder work.
$ws.Range("B8").Value = 4794
$ws.Range("C8").Value = "30/03/2020"
$ws.Range("D8").Value = 0.86805555555555547   # 8:50 PM
$ws.Range("E8").Value = 0.95138888888888884   # 10:50 PM

# --- Row 9 (new entry, 31/03/2020): finished the arithmetic unit.
$ws.Range("B9").Value = 4794
$ws.Range("C9").Value = "31/03/2020"
$ws.Range("D9").Value = 0.8618055555555556    # 8:41 PM
$ws.Range("E9").Value = 0.43055555555555558   # 10:20 AM

$ws.Range("G8").Value = "Completed full adder and ripple adder"
$ws.Range("G9").Value = "Finished Arithmetic unit"

# --- Active selection moves to G10, matching the new last-edited row.
[void]$ws.Range("G10").Select()
